$wb = $excel.ActiveWorkbook
$wsMF = $wb.Worksheets.Item("Mass_Fractions")
$wsUN = $wb.Worksheets.Item("Uncertainties")

# Corrected average mass fraction data
$rowVals = @([double]"2", [double]"7.7650977849269817E-2", [double]"2.295196542455219E-3", [double]"0.10744183029275578", [double]"7.4888766449407547E-3", [double]"2.2146447326846425E-2", [double]"0", [double]"0", [double]"0.46538946313854235", [double]"0.2343748911147219", [double]"6.178127487171072E-3", [double]"0", [double]"5.2032017840304118E-2", [double]"1.6932368035751037E-2", [double]"0", [double]"6.2143403694252177E-4", [double]"4.8065731868859288E-4", [double]"6.9677123716103276E-3")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(2, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"4", [double]"8.8263440384111477E-2", [double]"2.8179217619136339E-3", [double]"0.12255206607055497", [double]"8.9394414153531696E-3", [double]"2.6900451809774916E-2", [double]"0", [double]"0", [double]"0.5598731259137949", [double]"9.6301838351008376E-2", [double]"9.2042520986496453E-3", [double]"0", [double]"5.9865686843953922E-2", [double]"1.6670157757145233E-2", [double]"0", [double]"7.1433700159036532E-4", [double]"4.2827958251172407E-4", [double]"7.4690010096376278E-3")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(3, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"6", [double]"9.5287438272993782E-2", [double]"2.5434322278413018E-3", [double]"0.1154708949988751", [double]"9.7726810599756964E-3", [double]"3.1601634040896359E-2", [double]"0", [double]"0", [double]"0.6095917336643093", [double]"4.363583023651358E-2", [double]"9.0725149215348019E-3", [double]"0", [double]"5.9615535336234364E-2", [double]"1.553961581292968E-2", [double]"0", [double]"6.7114124664696196E-4", [double]"3.7620800673437602E-4", [double]"6.8213401745146097E-3")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(4, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"10", [double]"9.1756517114560787E-2", [double]"1.7365967378778189E-3", [double]"0.10684336666692139", [double]"1.0713486432351168E-2", [double]"5.1182742407965803E-2", [double]"0", [double]"0", [double]"0.66058954418115745", [double]"1.1108883206032259E-2", [double]"6.2023408613692379E-3", [double]"0", [double]"4.5212093447173404E-2", [double]"9.4433982378554285E-3", [double]"0", [double]"5.3702381212209925E-4", [double]"2.3844864651348271E-4", [double]"4.4355582480997083E-3")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(5, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"14", [double]"7.9713082848292233E-2", [double]"1.4330728065518248E-3", [double]"8.4488738168766647E-2", [double]"1.1199782240325801E-2", [double]"7.1896651146980878E-2", [double]"0", [double]"0", [double]"0.69539153922957353", [double]"5.7761367918317376E-3", [double]"4.6019546562011794E-3", [double]"0", [double]"3.6160433844838344E-2", [double]"5.8589985020059765E-3", [double]"0", [double]"2.4419956388441769E-4", [double]"1.7734380374179957E-4", [double]"3.0580663970055195E-3")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(6, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"20", [double]"6.7253411155698606E-2", [double]"5.7863736106888563E-4", [double]"5.809569063464963E-2", [double]"1.170222505602906E-2", [double]"0.11096195447465579", [double]"0", [double]"0", [double]"0.72739077718590006", [double]"1.1239602370112481E-3", [double]"1.9153446624747194E-3", [double]"0", [double]"1.7193818265629542E-2", [double]"2.5360880986783027E-3", [double]"0", [double]"0", [double]"8.5107857982037993E-5", [double]"1.1629850102219822E-3")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(7, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"30", [double]"5.3929007797550818E-2", [double]"1.7033429176696452E-4", [double]"4.1518131044865604E-2", [double]"1.18935167730075E-2", [double]"0.14280719567999819", [double]"0", [double]"0", [double]"0.74197615175252507", [double]"1.3150727049635277E-4", [double]"5.1464200237591191E-4", [double]"0", [double]"5.9870608844036291E-3", [double]"7.9529589593393827E-4", [double]"0", [double]"0", [double]"1.4086821208677452E-5", [double]"2.6306978586730444E-4")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(8, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"45", [double]"3.3242976412738948E-2", [double]"0", [double]"2.1168619157900498E-2", [double]"1.1955650000754879E-2", [double]"0.17964657220773519", [double]"0", [double]"0", [double]"0.75351287490926955", [double]"9.0038145856967015E-6", [double]"2.1709072341432119E-5", [double]"0", [double]"3.4008612768205479E-4", [double]"0", [double]"0", [double]"0", [double]"0", [double]"1.0250829699170914E-4")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(9, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"60", [double]"2.0073044886322281E-2", [double]"0", [double]"1.2944034505505192E-2", [double]"1.2141731555582434E-2", [double]"0.20270824766078444", [double]"0", [double]"0", [double]"0.75210310613450371", [double]"0", [double]"1.0209353094297727E-5", [double]"0", [double]"1.9625904207698935E-5", [double]"0", [double]"0", [double]"0", [double]"0", [double]"0")
for ($c = 1; $c -le 18; $c++) { $wsMF.Cells.Item(10, $c).Value = $rowVals[$c-1] }

# Corrected uncertainty data
$rowVals = @([double]"2", [double]"6.0854384124294121E-3", [double]"2.03061726349037E-4", [double]"2.2023274547215602E-2", [double]"1.0255830932353865E-3", [double]"8.3302052743043069E-3", [double]"0", [double]"0", [double]"5.2654665694891058E-2", [double]"0.11228507448389466", [double]"5.6339607663276337E-4", [double]"0", [double]"4.7873504517861772E-3", [double]"3.2283518853831496E-3", [double]"0", [double]"2.8891284009774262E-4", [double]"4.6898734430133737E-5", [double]"1.0714781846085442E-3")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(2, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"4", [double]"8.8524209316939162E-3", [double]"5.5416065245914091E-4", [double]"3.3987479616525512E-2", [double]"9.5609334241487009E-4", [double]"2.1465242817737092E-2", [double]"0", [double]"0", [double]"4.8811032714784708E-2", [double]"3.9008751422410136E-2", [double]"2.6170880608307068E-3", [double]"0", [double]"1.2668045467365984E-2", [double]"4.4865863742205182E-3", [double]"0", [double]"1.5683845986274485E-4", [double]"1.3154169593135483E-4", [double]"2.0393836793288068E-3")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(3, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"6", [double]"4.4079630298698973E-3", [double]"1.5391509782393262E-4", [double]"3.6014999045295953E-2", [double]"7.8327552523534735E-4", [double]"1.6302779230136398E-3", [double]"0", [double]"0", [double]"2.7084795877457882E-2", [double]"8.497893378592164E-3", [double]"9.2272181881775952E-4", [double]"0", [double]"2.676127108047177E-3", [double]"1.5175312650972436E-3", [double]"0", [double]"1.8995409501206833E-4", [double]"2.9627595365873751E-5", [double]"4.6730123102439719E-4")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(4, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"10", [double]"4.5668406785444324E-3", [double]"4.436437291757716E-4", [double]"2.0798768914459604E-2", [double]"1.0523905286749257E-3", [double]"1.6443486514894828E-2", [double]"0", [double]"0", [double]"4.8538962792641076E-2", [double]"2.5520288486608131E-3", [double]"1.088596957443091E-3", [double]"0", [double]"7.7556596722893993E-3", [double]"6.8678432351284375E-4", [double]"0", [double]"1.4845374855898826E-4", [double]"4.0665997843380784E-5", [double]"8.4807760215003885E-4")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(5, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"14", [double]"1.2567529737204958E-2", [double]"1.1380323161611884E-3", [double]"1.8735108498394423E-2", [double]"1.0471437536660785E-3", [double]"4.8488105066337661E-2", [double]"0", [double]"0", [double]"4.7632013639611233E-2", [double]"7.3085470611360401E-3", [double]"3.5258331228958874E-3", [double]"0", [double]"2.4281521895849967E-2", [double]"3.2795987926104906E-3", [double]"0", [double]"1.4932676024025438E-4", [double]"3.3800210097645087E-5", [double]"2.4862407313451353E-3")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(6, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"20", [double]"1.6313339045567433E-3", [double]"2.765393325942855E-5", [double]"1.102171382077682E-2", [double]"8.2979112546681161E-4", [double]"2.3895743913142009E-3", [double]"0", [double]"0", [double]"1.5551085334641232E-2", [double]"1.002582755334319E-4", [double]"1.0104146687145279E-4", [double]"0", [double]"4.0929596616772824E-4", [double]"3.9322038916370988E-5", [double]"0", [double]"0", [double]"3.7503707853534392E-6", [double]"6.3430885460600777E-5")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(7, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"30", [double]"2.2191814466303371E-3", [double]"4.2330919110245834E-5", [double]"1.2527150507768722E-2", [double]"9.5366824485095566E-4", [double]"1.0797938089490413E-2", [double]"0", [double]"0", [double]"3.3998121747070432E-2", [double]"9.8677994649125857E-5", [double]"1.6096998566692458E-4", [double]"0", [double]"1.5135767352856199E-3", [double]"2.5997814127629693E-4", [double]"0", [double]"0", [double]"2.8190353332233043E-5", [double]"9.3734850346650212E-5")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(8, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"45", [double]"3.3655534399250511E-3", [double]"0", [double]"5.4989565239554053E-3", [double]"8.845308428185282E-4", [double]"8.9765025241268434E-3", [double]"0", [double]"0", [double]"2.4884511144939093E-2", [double]"6.5138348049137553E-6", [double]"4.8315652616327985E-6", [double]"0", [double]"6.8026281944652588E-4", [double]"0", [double]"0", [double]"0", [double]"0", [double]"2.0516726322287302E-4")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(9, $c).Value = $rowVals[$c-1] }
$rowVals = @([double]"60", [double]"1.8715828853326263E-3", [double]"0", [double]"2.8302195079080296E-3", [double]"9.8300279107332422E-4", [double]"9.5255871993857547E-3", [double]"0", [double]"0", [double]"2.35997853645207E-2", [double]"0", [double]"2.0433009624408887E-5", [double]"0", [double]"3.9257565520056565E-5", [double]"0", [double]"0", [double]"0", [double]"0", [double]"0")
for ($c = 1; $c -le 18; $c++) { $wsUN.Cells.Item(10, $c).Value = $rowVals[$c-1] }
